$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.771.56'
$ws.Range("E2").Value = '  +2.06%  '
$ws.Range("D3").Value = '1.876.09'
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.37'
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  +1.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3921'
$ws.Range("E8").Value = '  +1.74%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07922'
$ws.Range("E9").Value = '  +1.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9753'
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.35'
$ws.Range("E11").Value = '  +2.15%  '
$ws.Range("D12").Value = '1.841.60'
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.752'
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.956'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07003'
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.61'
$ws.Range("E16").Value = '  +2.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001012'
$ws.Range("E18").Value = '  +1.87%  '
$ws.Range("E19").Value = '  +0.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '28.773.90'
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.339'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  +1.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.114'
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("D25").Value = '2.073.41'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.63'
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.39'
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.761'
$ws.Range("E28").Value = '  +1.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.006'
$ws.Range("E29").Value = '  +1.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.70'
$ws.Range("E30").Value = '  +2.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09383'
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9400'
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.330'
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.353'
$ws.Range("E34").Value = '  +2.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.356'
$ws.Range("E35").Value = '  -2.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05870'
$ws.Range("E36").Value = '  -2.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02117'
$ws.Range("E37").Value = '  -1.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.148'
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.932'
$ws.Range("E39").Value = '  +4.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5677'
$ws.Range("E40").Value = '  +1.26%  '
$ws.Range("E41").Value = '  +1.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.965'
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07242'
$ws.Range("E43").Value = '  +3.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.84'
$ws.Range("E44").Value = '  +2.34%  '
$ws.Range("E45").Value = '  +1.11%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.137'
$ws.Range("E46").Value = '  -8.90%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.125'
$ws.Range("E47").Value = '  -5.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.853'
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '113.91'
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.360'
$ws.Range("E50").Value = '  +1.25%  '
$ws.Range("E51").Value = '  +0.22%  '
